$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2-7: rename labels with "_50" suffix and refresh stats ---

# Row 2 (a1 -> a1_50)
$ws.Range("B2").Value = "a1_50"
$ws.Range("L2").Value = 633.0726
$ws.Range("M2").Value = 10.55
$ws.Range("N2").Value = 1.341614484786987
$ws.Range("O2").Value = 33
$ws.Range("P2").Value = 400.4725646972656
$ws.Range("Q2").Value = 0.21
$ws.Range("R2").Value = 321.04

# Row 3 (a2 -> a2_50)
$ws.Range("B3").Value = "a2_50"
$ws.Range("L3").Value = 660.1042
$ws.Range("M3").Value = 11
$ws.Range("N3").Value = 1.133472442626953
$ws.Range("O3").Value = 41
$ws.Range("P3").Value = 206.3318634033203
$ws.Range("Q3").Value = 0.25
$ws.Range("R3").Value = 285.16

# Row 4 (a3 -> a3_50)
$ws.Range("B4").Value = "a3_50"
$ws.Range("L4").Value = 636.4971
$ws.Range("M4").Value = 10.61
$ws.Range("N4").Value = 1.348449468612671
$ws.Range("O4").Value = 36
$ws.Range("P4").Value = 167.7512664794922
$ws.Range("Q4").Value = 0.25
$ws.Range("R4").Value = 319.93

# Row 5 (b -> b_50)
$ws.Range("B5").Value = "b_50"
$ws.Range("L5").Value = 636.4864
$ws.Range("M5").Value = 10.61
$ws.Range("N5").Value = 1.349513411521912
$ws.Range("O5").Value = 36
$ws.Range("P5").Value = 254.3524932861328
$ws.Range("Q5").Value = 0.4
$ws.Range("R5").Value = 321.09

# Row 6 (c -> c_50)
$ws.Range("B6").Value = "c_50"
$ws.Range("L6").Value = 669.2868
$ws.Range("M6").Value = 11.15
$ws.Range("N6").Value = 1.340455532073975
$ws.Range("O6").Value = 36
$ws.Range("P6").Value = 163.3970031738281
$ws.Range("Q6").Value = 0.28
$ws.Range("R6").Value = 320.76

# Row 7 (d -> d_50)
$ws.Range("B7").Value = "d_50"
$ws.Range("L7").Value = 655.01
$ws.Range("M7").Value = 10.92
$ws.Range("Q7").Value = 0.34
$ws.Range("R7").Value = 278.31

# --- Append new rows 8-13: fresh runs without the "_50" suffix (epochs=20) ---

$newRows = @(
    @{ A=6;  B="a1"; L=265.258;   M=4.42; N=1.345364928245544; O=20; P=137.91455078125;    Q=0.45; R=321.32 },
    @{ A=7;  B="a2"; L=260.7209;  M=4.35; N=1.202785968780518; O=20; P=120.4608306884766;   Q=0.24; R=280.73 },
    @{ A=8;  B="a3"; L=262.6342;  M=4.38; N=1.377196192741394; O=20; P=177.4877624511719;   Q=0.35; R=321.58 },
    @{ A=9;  B="b";  L=265.6575;  M=4.43; N=1.343576550483704; O=20; P=188.3482666015625;   Q=0.27; R=322.53 },
    @{ A=10; B="c";  L=267.9921;  M=4.47; N=1.347642779350281; O=20; P=138.1145477294922;   Q=0.26; R=321.44 },
    @{ A=11; B="d";  L=255.7635;  M=4.26; N=1.210534691810608; O=20; P=399.2044982910156;   Q=0.33; R=285.27 }
)

$ws.Range("A7").Copy()

$rowIndex = 8
foreach ($rec in $newRows) {
    $ws.Cells.Item($rowIndex, 1).PasteSpecial(-4122)
    $ws.Cells.Item($rowIndex, 1).Value = $rec.A
    $ws.Cells.Item($rowIndex, 2).Value = $rec.B
    $ws.Cells.Item($rowIndex, 3).Value = 1024
    $ws.Cells.Item($rowIndex, 4).Value = 256
    $ws.Cells.Item($rowIndex, 5).Value = 4
    $ws.Cells.Item($rowIndex, 6).Value = 6
    $ws.Cells.Item($rowIndex, 7).Value = 0
    $ws.Cells.Item($rowIndex, 8).Value = 0.001
    $ws.Cells.Item($rowIndex, 9).Value = 20
    $ws.Cells.Item($rowIndex, 10).Value = 4
    $ws.Cells.Item($rowIndex, 11).Value = "yes"
    $ws.Cells.Item($rowIndex, 12).Value = $rec.L
    $ws.Cells.Item($rowIndex, 13).Value = $rec.M
    $ws.Cells.Item($rowIndex, 14).Value = $rec.N
    $ws.Cells.Item($rowIndex, 15).Value = $rec.O
    $ws.Cells.Item($rowIndex, 16).Value = $rec.P
    $ws.Cells.Item($rowIndex, 17).Value = $rec.Q
    $ws.Cells.Item($rowIndex, 18).Value = $rec.R
    $ws.Cells.Item($rowIndex, 19).Value = 1
    $rowIndex = $rowIndex + 1
}
